$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '65.164.35', '  +0.47%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '2.946.47', '  -0.98%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  -0.12%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '567.85', '  -2.31%  '),
    @(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '158.70', '  +3.54%  '),
    @(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  +0.00%  '),
    @(8, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.519', '  +1.06%  '),
    @(9, 'LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '2.941.62', '  -1.05%  '),
    @(10, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '6.70', '  -3.69%  '),
    @(11, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.151', '  +0.70%  '),
    @(12, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.457', '  +2.38%  '),
    @(13, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000245', '  +3.35%  '),
    @(14, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '34.10', '  +0.75%  '),
    @(15, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.126', '  -0.49%  '),
    @(16, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '65.356.90', '  +0.70%  '),
    @(17, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '3.436.90', '  -1.13%  '),
    @(18, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '6.96', '  +0.97%  '),
    @(19, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '2.949.60', '  -0.97%  '),
    @(20, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '446.45', '  -0.43%  '),
    @(21, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '13.83', '  +1.45%  '),
    @(22, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.681', '  +0.32%  '),
    @(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '7.24', '  -0.62%  '),
    @(24, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '82.86', '  +2.40%  '),
    @(25, 'Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '2.18', '  -0.36%  '),
    @(26, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '12.05', '  -2.35%  '),
    @(27, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.999', '  -0.06%  '),
    @(28, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '9.97', '  -6.06%  '),
    @(29, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '7.95', '  +1.54%  '),
    @(30, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '2.34', '  -2.21%  '),
    @(31, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '2.57', '  -0.57%  '),
    @(32, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0₃0986', '  -3.01%  '),
    @(33, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '27.22', '  +1.99%  '),
    @(34, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.109', '  -0.64%  '),
    @(35, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '1.00', '  -0.12%  '),
    @(36, 'Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '0.974', '  -0.68%  '),
    @(37, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '5.72', '  +0.39%  '),
    @(38, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '49.04', '  +0.34%  '),
    @(39, 'Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '1.98', '  -4.65%  '),
    @(40, 'Arweave', 'https://coinranking.com/coin/7XWg41D1+arweave-ar', '43.32', '  -3.17%  '),
    @(41, 'TheGraph', 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt', '0.298', '  -0.12%  '),
    @(42, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.119', '  -0.91%  '),
    @(43, 'dogwifhat', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif', '2.78', '  -4.10%  '),
    @(44, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.42', '  +0.44%  '),
    @(45, 'Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '383.54', '  +0.55%  '),
    @(46, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0351', '  +0.80%  '),
    @(47, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '2.738.58', '  -1.02%  '),
    @(48, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '132.12', '  -1.79%  '),
    @(49, 'USDe', 'https://coinranking.com/coin/exbfr2U-0+usde-usde', '1.00', '  +0.02%  '),
    @(50, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.107', '  +1.42%  '),
    @(51, 'ThetaToken', 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta', '2.15', '  +5.58%  '),
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).NumberFormat = '@'
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
}

